# Update NATMI LR-pair TPM-derived statistics on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5909176666666667
$ws.Range("H2").Value = 1.772753
$ws.Range("I2").Value = 0.9937758428931484
$ws.Range("J2").Value = 0.9937758428931482
$ws.Range("P2").Value = 0.9810128591839974
$ws.Range("Q2").Value = 0.006085861049000001
$ws.Range("R2").Value = 0.054772749441
$ws.Range("S2").Value = 0.9749068810245946
$ws.Range("T2").Value = 0.9749068810245943

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5909176666666667
$ws.Range("H3").Value = 1.772753
$ws.Range("I3").Value = 0.9937758428931484
$ws.Range("J3").Value = 0.9937758428931482
$ws.Range("Q3").Value = 0.0001177895882222222
$ws.Range("R3").Value = 0.001060106294
$ws.Range("S3").Value = 0.01886896186855382
$ws.Range("T3").Value = 0.01886896186855382

# Row 4
$ws.Range("I4").Value = 0.006224157106851674
$ws.Range("J4").Value = 0.006224157106851673
$ws.Range("P4").Value = 0.9810128591839974
$ws.Range("R4").Value = 0.000343049391
$ws.Range("S4").Value = 0.006105978159402958
$ws.Range("T4").Value = 0.006105978159402958

# Row 5
$ws.Range("I5").Value = 0.006224157106851674
$ws.Range("J5").Value = 0.006224157106851673
$ws.Range("S5").Value = 0.0001181789474487157
$ws.Range("T5").Value = 0.0001181789474487157
